# Weekly update: insert a new record at row 216 (pushing existing rows 216:325
# down to 217:326) and populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 216; this shifts rows
# 216:325 down to 217:326 and grows the sheet to A1:R326.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new weekly record.
$ws.Cells.Item(216, 1).Value = 6
$ws.Cells.Item(216, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(216, 3).Value = "Metropolitana"
$ws.Cells.Item(216, 4).Value = 45097
$ws.Cells.Item(216, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 5).Value = 13
$ws.Cells.Item(216, 6).Value = 100112001
$ws.Cells.Item(216, 7).Value = "Berenjena"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 290
$ws.Cells.Item(216, 11).Value = 5000
$ws.Cells.Item(216, 12).Value = 6000
$ws.Cells.Item(216, 13).Value = 5586
$ws.Cells.Item(216, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(216, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(216, 16).Value = 112
$ws.Cells.Item(216, 17).Value = 50
$ws.Cells.Item(216, 18).Value = "Hortaliza"
